# "maj template comment à la fin"
# Move the "Comment" column (currently column J, together with its
# header/type/format/example rows) to the end of the header block
# (column N), shifting the SamplePreparationDate/Androstenone/Scatol/
# Indole block (columns K:N) one column to the left (J:M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 1..5) {
    $commentCol = $ws.Cells.Item($row, 10).Value()

    $ws.Cells.Item($row, 10).Value = $ws.Cells.Item($row, 11).Value()
    $ws.Cells.Item($row, 11).Value = $ws.Cells.Item($row, 12).Value()
    $ws.Cells.Item($row, 12).Value = $ws.Cells.Item($row, 13).Value()
    $ws.Cells.Item($row, 13).Value = $ws.Cells.Item($row, 14).Value()

    $ws.Cells.Item($row, 14).Value = $commentCol
}
